$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.526.57'
$ws.Range('D3').Value = '2.469.97'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.79'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.93'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.38'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '2.849.39'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.84'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.00'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.40%  '
$ws.Range('D16').Value = '2.465.04'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').Value = '41.559.28'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.51'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.15'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.73'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.38'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.87'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.08%  '
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0757'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.25'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.37%  '
$ws.Range('E36').Value = '  -6.77%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.103'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.114'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('E39').Value = '  -5.70%  '
$ws.Range('E40').Value = '  -10.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.01'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.54%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').Value = '1.943.56'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.55'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.06'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('D48').Value = '2.707.44'
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.81'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.25'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.60%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.172'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.63%  '
